# Add new test data to the "nested" sheet (sheet2):
#   - F1 gets a new label "aaa.bbb" (new shared string)
#   - A7 gets the same label "aaa.bbb"
#   - Selection moves to A8 (next empty cell below the new data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "aaa.bbb"
$ws.Range("A7").Value = "aaa.bbb"

$ws.Range("A8").Select()
